$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @(date, D, E, F, G, H) reflecting the updated attendance data
$data = @{
    3  = @("28-07-2022", 1, 0, 0, 1, 1)
    4  = @("01-08-2022", 1, 1, 0, 0, 0)
    5  = @("04-08-2022", 1, 1, 0, 0, 0)
    6  = @("08-08-2022", 1, 1, 0, 0, 0)
    7  = @("11-08-2022", 0, 0, 0, 0, 1)
    8  = @("15-08-2022", 0, 0, 0, 0, 1)
    9  = @("18-08-2022", 0, 0, 0, 0, 1)
    10 = @("22-08-2022", 1, 1, 0, 0, 0)
    11 = @("25-08-2022", 1, 1, 0, 0, 0)
    12 = @("29-08-2022", 1, 1, 0, 0, 0)
    13 = @("01-09-2022", 1, 1, 0, 0, 0)
    14 = @("05-09-2022", 0, 0, 0, 0, 1)
    15 = @("08-09-2022", 1, 1, 0, 0, 0)
    16 = @("12-09-2022", 1, 1, 0, 0, 0)
    17 = @("15-09-2022", 0, 0, 0, 0, 1)
    18 = @("19-09-2022", 0, 0, 0, 0, 1)
    19 = @("22-09-2022", 0, 0, 0, 0, 1)
    20 = @("26-09-2022", 1, 1, 0, 0, 0)
    21 = @("29-09-2022", 0, 0, 0, 0, 1)
}

# Rows whose DD-MM-YYYY text is ambiguous with MM-DD-YYYY (day <= 12) need to be
# forced to Text format first, otherwise Excel auto-converts the literal string
# into a date serial number on assignment.
$needsTextFormat = @(4, 5, 6, 7, 13, 14, 15, 16)

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $cellA = $ws.Cells.Item($row, 1)
    if ($needsTextFormat -contains $row) {
        $cellA.NumberFormat = "@"
    }
    $cellA.Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
    $ws.Cells.Item($row, 8).Value = $vals[5]
}
